$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "_old" / "_new" header-column suffixes to the respective
#    format-version names ("_FV2410" resp. "_FV2504"), as described by the
#    commit "Use <formatversion> as suffix for table headers".
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2410"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2504"
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into a real Excel Table ("Table1") so the new
#    header names also show up as the table's column headers.
#
#    The header row (A1:U1) already carries manual formatting (bold font,
#    grey fill, border, centered+wrapped alignment). If a ListObject is
#    created while that formatting is already present, Excel captures it as
#    a one-off "headerRowDxfId" override (and adds a new <dxf> to
#    styles.xml). To keep styles.xml untouched, the header formatting is
#    stashed away, cleared, the table is added against plain cells, and the
#    original formatting is then re-applied with a simple format copy/paste
#    (which Excel stores as normal cell styles, not as a table dxf).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stashRange = $ws.Range("A67:U67")

$stashRange.Value = $headerRange.Value
$headerRange.Copy()
$stashRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)  # xlSrcRange, xlYes
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$stashRange.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, top-left cell of the
#    scrollable area is A2).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
